$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("A4").Value = 111633837
$ws.Range("B4").Value = 98535
$ws.Range("D4").Value = "LC"
$ws.Range("E4").Value = 222498
$ws.Range("F4").Value = "Blåsippa"
$ws.Range("G4").Value = "Hepatica nobilis"
$ws.Range("H4").Value = "Schreb."

# Row 5
$ws.Range("A5").Value = 111633890
$ws.Range("B5").Value = 90658
$ws.Range("D5").Value = "NT"
$ws.Range("E5").Value = 4361
$ws.Range("F5").Value = "Orange taggsvamp"
$ws.Range("G5").Value = "Hydnellum aurantiacum"
$ws.Range("H5").Value = "(Batsch:Fr.) P.Karst."

# Row 6
$ws.Range("A6").Value = 111633843
$ws.Range("Q6").Value = 676486.710397501
$ws.Range("R6").Value = 6618439.724061669

# Row 7
$ws.Range("A7").Value = 111634290
$ws.Range("Q7").Value = 676708.8668162767
$ws.Range("R7").Value = 6618511.450801066

# Row 8
$ws.Range("A8").Value = 111634304
$ws.Range("B8").Value = 90687
$ws.Range("E8").Value = 5964
$ws.Range("F8").Value = "Fjällig taggsvamp s.str."
$ws.Range("G8").Value = "Sarcodon imbricatus s.str."
$ws.Range("H8").Value = "(L.:Fr.) P.Karst."
